# "Added top emails in config"
#
# Inserts a new configuration row ("topEmails" / 100 / description) above
# the existing "usedFolder" row on the Constants sheet, pushing every row
# from the old row 21 onward down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# Shift row 21 (and everything below it) down by inserting a blank row.
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row with the new constant.
$ws.Cells.Item(21, 1).Value = "topEmails"
$ws.Cells.Item(21, 2).Value = 100
$ws.Cells.Item(21, 3).Value = "How many emails to be processed at a time(only use integer numbers)"

# Match the author's final selection/viewport on the sheet.
$ws.Range("C22").Select()
